# "added logic for finding headers"
# The customer_id / customer_name / customer_password columns are
# reordered so that the numeric id column moves from first to last:
#   old: A=customer_id(number) B=customer_name C=customer_password
#   new: A=customer_name       B=customer_password C=customer_id(number)
# A new (blank, but formatted) row is also appended below the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "customer_name"
$ws.Range("B1").Value = "customer_password"
$ws.Range("C1").Value = "customer_id"

# --- Data rows ----------------------------------------------------------
# row, name, password, id
$rows = @(
    @(2, "shiva",             "shiva@437",   1),
    @(3, "swamy",              "swamy@123",  2),
    @(4, "sasanala grandham",  "bhavani@123",3),
    @(5, "sharuna",            "sharuna@123",4)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}

# --- New trailing (blank but formatted) row ------------------------------
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A6").ClearContents()
